$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing formulas in row 2 to use absolute row reference for K2
$ws.Range("E2").Formula = "=C2+150-K`$2"
$ws.Range("F2").Formula = "=C2+150+0.5*K`$2"
$ws.Range("G2").Formula = "=D2+300-K`$2"
$ws.Range("H2").Formula = "=D2+300-0.25*K`$2"

# Add new row 3 with data for a new location (tree graph / first_map)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "castle"
$ws.Range("C3").Value = 1024
$ws.Range("D3").Value = 896
$ws.Range("E3").Formula = "=C3+150-K`$2"
$ws.Range("F3").Formula = "=C3+150+0.5*K`$2"
$ws.Range("G3").Formula = "=D3+300-K`$2"
$ws.Range("H3").Formula = "=D3+300-0.25*K`$2"
$ws.Range("I3").Value = "castle.png"
$ws.Range("J3").Value = "first_map"

# Update selection to match target state
$ws.Range("C10").Select()
